# update A12 and A19
$wb = $excel.ActiveWorkbook

# --- A12: replace the TxHash / NFTID pair in row 2 ---
$wsA12 = $wb.Worksheets.Item("A12")
$wsA12.Range("A2").Value = "ibc/5264E6AB7F094942F58C755FAD71F7C5F7DC8F9EA6A58AA16A9BA9694678698D"
$wsA12.Range("B2").Value = "ark151"

# --- A19: replace the six ClassID values in column A (rows 2-7) ---
$wsA19 = $wb.Worksheets.Item("A19")
$wsA19.Range("A2").Value = "258107172329D10A2879B454820C1B0AB57E8124AFE49799A844A82769D833BB"
$wsA19.Range("A3").Value = "1A15BF46DE2706D09F7951267A986890377CB56857D4409188E857146474001D"
$wsA19.Range("A4").Value = "176A6C67BA6D9921028FE9B6C12599E1635C48F9DDFD86449A3E005D388ECD19"
$wsA19.Range("A5").Value = "17217F017ABD2F25F846B95D7E14DCFE535B9DE1801246B27B3B776E364376A6"
$wsA19.Range("A6").Value = "28DD99AE253AE7DB85820CD66C6628907A94452919258F1BDB97B5BB272A3399"
$wsA19.Range("A7").Value = "10CF3A5C81FC12B55189500E3506F3039B967625B94C70DB763720574CA7AE79"

# --- add a new "B7" sheet, cloned from "B6", placed right after it ---
$wsB6 = $wb.Worksheets.Item("B6")
$wsB6.Copy($null, $wsB6)
$wsB7 = $wb.Worksheets.Item("B6 (2)")
$wsB7.Name = "B7"

# --- move the active tab from A20 to A19, with selection K19 ---
$wsA19.Activate()
[void]$wsA19.Range("K19").Select()
